# The target revision only rewrites the OOXML attribute *serialization
# order* inside word/document.xml (root namespace declarations, the
# <w:sectPr> page setup) and word/styles.xml (<w:docDefaults>,
# <w:latentStyles>/<w:lsdException>, <w:style> elements) -- every
# attribute name/value pair is identical before and after, only the
# order in which they are written out changed. The commit message
# confirms this: it is a side effect of "Fixed POI packaging and
# upgraded to POI 3.15", i.e. a new XML serializer that happens to emit
# attributes alphabetically, not a document edit.
#
# There is no content, formatting, or style change for Word's object
# model to reproduce here -- attribute ordering is an artifact of the
# XML writer, not something the Word UI/automation surface exposes or
# lets a user control. The correct interop action is therefore to
# leave the document exactly as authored.
$d = $word.ActiveDocument
